$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Moorings")
$ws2 = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Update instrument / glider reference designators from GL002 to GL362 ---
# Moorings sheet: Ref Des (A2) and Deployment Number (C2)
$ws1.Range("A2").Value = "GP05MOAS-GL362"
$ws1.Range("C2").Value = 1

# Asset_Cal_Info sheet: Ref Des (A3:A9) and Deployment Number (C3:C9)
$ws2.Range("A3").Value = "GP05MOAS-GL362-00-ENG000000"
$ws2.Range("C3").Value = 1

$ws2.Range("A4").Value = "GP05MOAS-GL362-01-FLORDM000"
$ws2.Range("C4").Value = 1

$ws2.Range("A5").Value = "GP05MOAS-GL362-01-FLORDM000"
$ws2.Range("C5").Value = 1

$ws2.Range("A6").Value = "GP05MOAS-GL362-01-FLORDM000"
$ws2.Range("C6").Value = 1

$ws2.Range("A7").Value = "GP05MOAS-GL362-01-FLORDM000"
$ws2.Range("C7").Value = 1

$ws2.Range("A8").Value = "GP05MOAS-GL362-02-DOSTAM000"
$ws2.Range("C8").Value = 1

$ws2.Range("A9").Value = "GP05MOAS-GL362-04-CTDGVM000"
$ws2.Range("C9").Value = 1

# --- Update the active selection on each sheet ---
$ws1.Activate()
$ws1.Range("D22").Select()

$ws2.Activate()
$ws2.Range("C10").Select()
